$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Lab marks for the student (row 6)
$ws.Range("B6").Value = 10
$ws.Range("C6").Value = 7.5
$ws.Range("D6").Value = 9
$ws.Range("E6").Value = 10
$ws.Range("F6").Value = 9.5

# Feedback comments per lab column (row 7, merged B7:B14 ... F7:F14)
# Order of entry matches the original authoring order (controls shared-string table order)
$ws.Range("B7").Value = "Good start but next time please use the course template you can download from the virtual campus"
$ws.Range("D7").Value = "Please, use the correct packages for the sessions in the future. Tromino numbers are better consecutive to understand what is happening"
$ws.Range("E7").Value = "Please, use PDFs"
$ws.Range("C7").Value = "You need to use the template that can be downloaded from the virtual campus. Insertion: times in the table seems OK (the explanation said that it didn't). Bubble does not work as expected"
$ws.Range("F7").Value = "Very good but times were not very representative"

# Update the view: scrolled right one column, selection on F7:F14
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("F7:F14").Select()
